$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency price/volume snapshot (Price col D, Volume(1h) col E).
# Price column is stored as plain display text (e.g. "1.000", "30.303.90"), so
# force text format before assigning to stop Excel from re-typing it as a number.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "30.303.90"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.10%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.928.74"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -0.65%  "

$ws.Range("E4").Value = "  +0.01%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "249.22"
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.7175"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.69%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -4.90%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "27.82"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -3.04%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.07083"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -4.96%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.7885"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -3.53%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.07986"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -1.92%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.929.48"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -0.58%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.371"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -2.69%  "

$ws.Range("E15").Value = "  -0.65%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "14.63"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.72%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "30.299.85"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.19%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "257.03"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.97%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.000008093"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -2.78%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "5.754"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -2.16%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "2.183.13"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.48%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.9999"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "6.821"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.97%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "9.531"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -3.45%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "164.49"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +2.38%  "

$ws.Range("E27").Value = "  -1.73%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.265"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -7.34%  "

$ws.Range("E29").Value = "  -3.83%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.355"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.91%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.530"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -2.65%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.404"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -1.92%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.139"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -2.39%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.05126"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.52%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.268"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.62%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7447"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.58%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.772"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +1.19%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01980"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.50%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.796"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.33%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "78.01"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -4.57%  "

$ws.Range("E41").Value = "  -2.81%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.4498"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.96%  "

$ws.Range("E43").Value = "  -1.67%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.8448"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.49%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.9997"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "100.83"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -2.14%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "9.810"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -0.74%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "7.459"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.19%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "36.85"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.61%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "950.85"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +8.00%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.4207"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.07%  "
